# Apply cell updates described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.111.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.558.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.18%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.000'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '292.29'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3961'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3238'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.72%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.38'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07335'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.086'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.05'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.88%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.677'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001146'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.668'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.557.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06597'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9995'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.317'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.123.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.345'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.439'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.873'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.730.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.015'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.76%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.735'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08389'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.643'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -13.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.081'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.04%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02283'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.24%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06154'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.151'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.217'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2067'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.73%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.79'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5881'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.769'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5629'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.905'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.143'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06864'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.95%  '
